$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = '[0.5421730213202043, 13.10491909383067]'
$ws.Range("N2").Value = 0.03390699391187768
$ws.Range("O2").Value = 0.03390699391187768
$ws.Range("U2").Value = '[4.765362072965749, 12.210158371411785]'
$ws.Range("V2").Value = 0.000035240466822905602839455241
$ws.Range("W2").Value = 0.000035240466822905602839455241
$ws.Range("M3").Value = '[0.010535859336384235, 14.109867382285383]'
$ws.Range("N3").Value = 0.04967287044907787
$ws.Range("O3").Value = 0.04967287044907787
$ws.Range("Q3").Value = '[-2.7673689040630816, -0.1006315965113842]'
$ws.Range("R3").Value = 0.03563895790123506
$ws.Range("S3").Value = 0.03563895790123506
$ws.Range("U3").Value = '[5.429395286749703, 13.51606046818419]'
$ws.Range("V3").Value = 0.000023327955005036930038195622
$ws.Range("W3").Value = 0.000023327955005036930038195622
$ws.Range("Y3").Value = 0.4040840840840918
$ws.Range("Z3").Value = 11.11231231231254
$ws.Range("M4").Value = '[-0.35937052420678484, 13.960616806717844]'
$ws.Range("N4").Value = 0.06212244149622315
$ws.Range("O4").Value = 0.06212244149622315
$ws.Range("Q4").Value = '[-3.182474239672543, -0.2641579408423853]'
$ws.Range("R4").Value = 0.02167664886900855
$ws.Range("S4").Value = 0.02167664886900855
$ws.Range("U4").Value = '[5.25556013683389, 13.325136241786005]'
$ws.Range("V4").Value = 0.000030420907374928418676063302
$ws.Range("W4").Value = 0.000030420907374928418676063302
$ws.Range("Y4").Value = 1.060720720720743
$ws.Range("Z4").Value = 12.77915915915941
$ws.Range("M5").Value = '[-1.031846092105722, 14.471793197137139]'
$ws.Range("N5").Value = 0.08763431059698479
$ws.Range("O5").Value = 0.08763431059698479
$ws.Range("Q5").Value = '[-3.849158566560468, -0.2012631930227693]'
$ws.Range("R5").Value = 0.03033121292927854
$ws.Range("S5").Value = 0.03033121292927854
$ws.Range("U5").Value = '[5.148674261202096, 13.302442519593761]'
$ws.Range("V5").Value = 0.00003947087820010963810091198
$ws.Range("W5").Value = 0.00003947087820010963810091198
$ws.Range("Y5").Value = 0.8081681681681845
$ws.Range("Z5").Value = 15.45621621621653
$ws.Range("M6").Value = '[-1.2270646350737229, 14.793325552632254]'
$ws.Range("N6").Value = 0.09498328970761771
$ws.Range("O6").Value = 0.09498328970761771
$ws.Range("Q6").Value = '[-5.4466851611787, 0.8176317216550011]'
$ws.Range("R6").Value = 0.1436379007227098
$ws.Range("S6").Value = 0.1436379007227098
$ws.Range("U6").Value = '[5.12589877850044, 13.347273039764554]'
$ws.Range("V6").Value = 0.000043805183582712459156027762
$ws.Range("W6").Value = 0.000043805183582712459156027762
$ws.Range("Y6").Value = -3.283183183183249
$ws.Range("Z6").Value = 21.87105105105149
$ws.Range("M7").Value = '[-1.256895988949096, 14.866122968479544]'
$ws.Range("N7").Value = 0.09601912210720487
$ws.Range("O7").Value = 0.09601912210720487
$ws.Range("Q7").Value = '[-4.729685036035086, -0.4780000834290763]'
$ws.Range("R7").Value = 0.01749774854600372
$ws.Range("S7").Value = 0.01749774854600372
$ws.Range("U7").Value = '[5.0697596541581085, 13.229507901470331]'
$ws.Range("V7").Value = 0.000045065775079011771708792367
$ws.Range("W7").Value = 0.000045065775079011771708792367
$ws.Range("Y7").Value = 1.919399399399431
$ws.Range("Z7").Value = 18.99195195195234
$ws.Range("M8").Value = '[-0.38666915139748603, 14.024640120866275]'
$ws.Range("N8").Value = 0.06304660226229042
$ws.Range("O8").Value = 0.06304660226229042
$ws.Range("Q8").Value = '[-4.503263943884468, -1.257894956392311]'
$ws.Range("R8").Value = 0.0008493165035166061
$ws.Range("S8").Value = 0.0008493165035166061
$ws.Range("U8").Value = '[5.052781712019428, 13.168120288123607]'
$ws.Range("V8").Value = 0.000044302189695866189822481829
$ws.Range("W8").Value = 0.000044302189695866189822481829
$ws.Range("Y8").Value = 5.051051051051155
$ws.Range("Z8").Value = 18.08276276276312
$ws.Range("M9").Value = '[-0.7998486858797236, 14.58595317458538]'
$ws.Range("N9").Value = 0.07781726301977621
$ws.Range("O9").Value = 0.07781726301977621
$ws.Range("Q9").Value = '[0.8490790955648087, 4.446658670846814]'
$ws.Range("R9").Value = 0.004830425816805484
$ws.Range("S9").Value = 0.004830425816805484
$ws.Range("U9").Value = '[5.0772436763167885, 13.086871766015966]'
$ws.Range("V9").Value = 0.000038226302600641297476085495
$ws.Range("W9").Value = 0.000038226302600641297476085495
$ws.Range("Y9").Value = 6.991631631631724
$ws.Range("Z9").Value = 20.68756756756783
$ws.Range("M10").Value = '[-0.573672721254642, 14.145636237282066]'
$ws.Range("N10").Value = 0.06984817430381196
$ws.Range("O10").Value = 0.06984817430381196
$ws.Range("Q10").Value = '[1.0000264903318858, 4.396342872591122]'
$ws.Range("R10").Value = 0.002518791213439231
$ws.Range("S10").Value = 0.002518791213439231
$ws.Range("U10").Value = '[5.042848423156576, 12.996810070214574]'
$ws.Range("V10").Value = 0.000038170202399490087304383817
$ws.Range("W10").Value = 0.000038170202399490087304383817
$ws.Range("Y10").Value = 7.183183183183272
$ws.Range("Z10").Value = 20.11291291291316
$ws.Range("M11").Value = '[-0.6366779821237074, 14.358725081082298]'
$ws.Range("N11").Value = 0.07191101117458487
$ws.Range("O11").Value = 0.07191101117458487
$ws.Range("Q11").Value = '[1.025184389459732, 4.421500771718968]'
$ws.Range("R11").Value = 0.00231510667046142
$ws.Range("S11").Value = 0.00231510667046142
$ws.Range("U11").Value = '[5.112196489310984, 13.095240458112341]'
$ws.Range("V11").Value = 0.000035107097127395903157776047
$ws.Range("W11").Value = 0.000035107097127395903157776047
$ws.Range("Y11").Value = 7.087407407407494
$ws.Range("Z11").Value = 20.01713713713739
$ws.Range("M12").Value = '[-0.3228455298972346, 13.931110655620156]'
$ws.Range("N12").Value = 0.06083931900375106
$ws.Range("O12").Value = 0.06083931900375106
$ws.Range("U12").Value = '[5.08085240530394, 13.043720022347127]'
$ws.Range("V12").Value = 0.00003618946456329653183327591
$ws.Range("W12").Value = 0.00003618946456329653183327591
$ws.Range("M13").Value = '[-0.43072155405570633, 13.981836161088674]'
$ws.Range("N13").Value = 0.06470177429913204
$ws.Range("O13").Value = 0.06470177429913204
$ws.Range("U13").Value = '[5.068079252210204, 13.05701663560831]'
$ws.Range("V13").Value = 0.000037979296126255952640877328
$ws.Range("W13").Value = 0.000037979296126255952640877328
$ws.Range("M14").Value = '[-0.4215070437718502, 14.028413801867867]'
$ws.Range("N14").Value = 0.06431102769347508
$ws.Range("O14").Value = 0.06431102769347508
$ws.Range("Q14").Value = '[1.3270791789938863, 4.34602707433543]'
$ws.Range("R14").Value = 0.0004523927187385368
$ws.Range("S14").Value = 0.0004523927187385368
$ws.Range("U14").Value = '[5.083298559295053, 13.042607262792114]'
$ws.Range("V14").Value = 0.000035909143268231247247967985
$ws.Range("W14").Value = 0.000035909143268231247247967985
$ws.Range("Y14").Value = 7.374734734734825
$ws.Range("Z14").Value = 18.86782782782807
$ws.Range("M15").Value = '[-0.5531807677137355, 14.272426632239718]'
$ws.Range("N15").Value = 0.06888216327508001
$ws.Range("O15").Value = 0.06888216327508001
$ws.Range("U15").Value = '[5.111685418942639, 13.098500194295521]'
$ws.Range("V15").Value = 0.000035276331031397489823575597
$ws.Range("W15").Value = 0.000035276331031397489823575597
$ws.Range("M16").Value = '[-0.24558570020544757, 13.86731986448869]'
$ws.Range("N16").Value = 0.05816291693807463
$ws.Range("O16").Value = 0.05816291693807463
$ws.Range("Q16").Value = '[1.4654476241970404, 4.358606023899352]'
$ws.Range("R16").Value = 0.0001968129344680491
$ws.Range("S16").Value = 0.0001968129344680491
$ws.Range("U16").Value = '[5.059418321762724, 13.025957579680021]'
$ws.Range("V16").Value = 0.000037634138853048781925281219
$ws.Range("W16").Value = 0.000037634138853048781925281219
$ws.Range("Y16").Value = 7.326846846846939
$ws.Range("Z16").Value = 18.34106106106129
$ws.Range("M17").Value = '[0.05097439504313428, 13.73774985082593]'
$ws.Range("N17").Value = 0.04838793180150502
$ws.Range("O17").Value = 0.04838793180150502
$ws.Range("Q17").Value = '[1.553500271144502, 4.396342872591122]'
$ws.Range("R17").Value = 0.0001185719985772327
$ws.Range("S17").Value = 0.0001185719985772327
$ws.Range("U17").Value = '[5.077581509481661, 13.050559120563321]'
$ws.Range("V17").Value = 0.000036772986935051932277929154
$ws.Range("W17").Value = 0.000036772986935051932277929154
$ws.Range("Y17").Value = 7.183183183183275
$ws.Range("Z17").Value = 18.00584584584607
